$wb = $excel.ActiveWorkbook

# ALC row 8 (Leve Item ID 4565)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 65.166664
$ws.Range("I8").Value = 65.166664
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 195.499992
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -56.49999199999999

# ALC row 9 (Leve Item ID 5487)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1227.2778
$ws.Range("I9").Value = 1320.125
$ws.Range("J9").Value = 484.5
$ws.Range("K9").Value = 1320.125
$ws.Range("L9").Value = 484.5
$ws.Range("M9").Value = -1151.125
$ws.Range("N9").Value = -822.5

# ALC row 103 (Leve Item ID 19909)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 810.2
$ws.Range("I103").Value = 299
$ws.Range("J103").Value = 1151
$ws.Range("K103").Value = 897
$ws.Range("L103").Value = 3453
$ws.Range("M103").Value = -311
$ws.Range("N103").Value = -4625

# ALC row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2522.5217
$ws.Range("I137").Value = 1651.5625
$ws.Range("J137").Value = 4513.2856
$ws.Range("K137").Value = 4954.6875
$ws.Range("L137").Value = 13539.8568
$ws.Range("M137").Value = -2404.6875
$ws.Range("N137").Value = -18639.8568

# ARM row 2 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5461.9585
$ws.Range("I2").Value = 5198.1113
$ws.Range("J2").Value = 6253.5
$ws.Range("K2").Value = 5198.1113
$ws.Range("L2").Value = 6253.5
$ws.Range("M2").Value = -5085.1113
$ws.Range("N2").Value = -6479.5

# ARM row 5 (Leve Item ID 5091)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 97.5
$ws.Range("I5").Value = 97.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 97.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 14.5

# ARM row 25 (Leve Item ID 2471)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 5000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 5000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -4598

# ARM row 45 (Leve Item ID 27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1995.05
$ws.Range("I45").Value = 1889.875
$ws.Range("J45").Value = 2415.75
$ws.Range("K45").Value = 1889.875
$ws.Range("L45").Value = 2415.75
$ws.Range("M45").Value = -1512.875
$ws.Range("N45").Value = -3169.75

# ARM row 116 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 5461.9585
$ws.Range("I116").Value = 5198.1113
$ws.Range("J116").Value = 6253.5
$ws.Range("K116").Value = 5198.1113
$ws.Range("L116").Value = 6253.5
$ws.Range("M116").Value = -2904.1113
$ws.Range("N116").Value = -10841.5

# ARM row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 34028.17
$ws.Range("I132").Value = 1761.826
$ws.Range("J132").Value = 157715.83
$ws.Range("K132").Value = 5285.478
$ws.Range("L132").Value = 473147.49
$ws.Range("M132").Value = -2755.478
$ws.Range("N132").Value = -478207.49

# BSM row 3 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5461.9585
$ws.Range("I3").Value = 5198.1113
$ws.Range("J3").Value = 6253.5
$ws.Range("K3").Value = 5198.1113
$ws.Range("L3").Value = 6253.5
$ws.Range("M3").Value = -5084.1113
$ws.Range("N3").Value = -6481.5

# BSM row 4 (Leve Item ID 5091)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 97.5
$ws.Range("I4").Value = 97.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 97.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 17.5

# BSM row 134 (Leve Item ID 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7584.8696
$ws.Range("I134").Value = 3903.1875
$ws.Range("J134").Value = 16000.143
$ws.Range("K134").Value = 11709.5625
$ws.Range("L134").Value = 48000.429
$ws.Range("M134").Value = -9174.5625
$ws.Range("N134").Value = -53070.429

# CRP row 4 (Leve Item ID 3742)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("N4").Value = 0

# CRP row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1895.3077
$ws.Range("I31").Value = 884.8570999999999
$ws.Range("J31").Value = 3074.1667
$ws.Range("K31").Value = 884.8570999999999
$ws.Range("L31").Value = 3074.1667
$ws.Range("M31").Value = -589.8570999999999
$ws.Range("N31").Value = -3664.1667

# CRP row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1895.3077
$ws.Range("I34").Value = 884.8570999999999
$ws.Range("J34").Value = 3074.1667
$ws.Range("K34").Value = 884.8570999999999
$ws.Range("L34").Value = 3074.1667
$ws.Range("M34").Value = -682.8570999999999
$ws.Range("N34").Value = -3478.1667

# CRP row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8026.091
$ws.Range("I58").Value = 6166.3335
$ws.Range("J58").Value = 16395
$ws.Range("K58").Value = 6166.3335
$ws.Range("L58").Value = 16395
$ws.Range("M58").Value = -5963.3335
$ws.Range("N58").Value = -16801

# CRP row 62 (Leve Item ID 12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 11476
$ws.Range("I62").Value = 7927
$ws.Range("J62").Value = 17686.75
$ws.Range("K62").Value = 7927
$ws.Range("L62").Value = 17686.75
$ws.Range("M62").Value = -7303
$ws.Range("N62").Value = -18934.75

# CRP row 65 (Leve Item ID 12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 11476
$ws.Range("I65").Value = 7927
$ws.Range("J65").Value = 17686.75
$ws.Range("K65").Value = 39635
$ws.Range("L65").Value = 88433.75
$ws.Range("M65").Value = -36515
$ws.Range("N65").Value = -94673.75

# CRP row 122 (Leve Item ID 36196)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3241.3333
$ws.Range("I122").Value = 3199.75
$ws.Range("J122").Value = 3262.125
$ws.Range("K122").Value = 9599.25
$ws.Range("L122").Value = 9786.375
$ws.Range("M122").Value = -7149.25
$ws.Range("N122").Value = -14686.375

# CRP row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8026.091
$ws.Range("I136").Value = 6166.3335
$ws.Range("J136").Value = 16395
$ws.Range("K136").Value = 18499.0005
$ws.Range("L136").Value = 49185
$ws.Range("M136").Value = -15949.0005
$ws.Range("N136").Value = -54285

# CUL row 5 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 143839.66
$ws.Range("I5").Value = 915.7406999999999
$ws.Range("J5").Value = 626207.9
$ws.Range("K5").Value = 2747.2221
$ws.Range("L5").Value = 1878623.7
$ws.Range("M5").Value = -2635.2221
$ws.Range("N5").Value = -1878847.7

# CUL row 55 (Leve Item ID 4733)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 339831.66
$ws.Range("I55").Value = 9495
$ws.Range("J55").Value = 505000
$ws.Range("K55").Value = 28485
$ws.Range("L55").Value = 1515000
$ws.Range("M55").Value = -28308
$ws.Range("N55").Value = -1515354

# CUL row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 30954348
$ws.Range("I131").Value = 48485650
$ws.Range("J131").Value = 19610564
$ws.Range("K131").Value = 145456950
$ws.Range("L131").Value = 58831692
$ws.Range("M131").Value = -145451910
$ws.Range("N131").Value = -58841772

# CUL row 132 (Leve Item ID 43972)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2997.4
$ws.Range("I132").Value = 1198.75
$ws.Range("J132").Value = 4196.5
$ws.Range("K132").Value = 10788.75
$ws.Range("L132").Value = 37768.5
$ws.Range("M132").Value = -8258.75
$ws.Range("N132").Value = -42828.5

# CUL row 135 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 143839.66
$ws.Range("I135").Value = 915.7406999999999
$ws.Range("J135").Value = 626207.9
$ws.Range("K135").Value = 8241.666299999999
$ws.Range("L135").Value = 5635871.100000001
$ws.Range("M135").Value = -5706.666299999999
$ws.Range("N135").Value = -5640941.100000001

# GSM row 97 (Leve Item ID 19940)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1287.4
$ws.Range("I97").Value = 1201.9445
$ws.Range("J97").Value = 1507.1428
$ws.Range("K97").Value = 1201.9445
$ws.Range("L97").Value = 1507.1428
$ws.Range("M97").Value = -705.9445000000001
$ws.Range("N97").Value = -2499.1428

# GSM row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11406.5
$ws.Range("I122").Value = 13685.5
$ws.Range("J122").Value = 10039.1
$ws.Range("K122").Value = 41056.5
$ws.Range("L122").Value = 30117.3
$ws.Range("M122").Value = -38606.5
$ws.Range("N122").Value = -35017.3

# GSM row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4912.706
$ws.Range("I132").Value = 4385.923
$ws.Range("J132").Value = 6624.75
$ws.Range("K132").Value = 13157.769
$ws.Range("L132").Value = 19874.25
$ws.Range("M132").Value = -10627.769
$ws.Range("N132").Value = -24934.25

# LTW row 7 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6208.684
$ws.Range("I7").Value = 6166.5
$ws.Range("J7").Value = 6219.933
$ws.Range("K7").Value = 6166.5
$ws.Range("L7").Value = 6219.933
$ws.Range("M7").Value = -6054.5
$ws.Range("N7").Value = -6443.933

# LTW row 31 (Leve Item ID 3043)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 10004614
$ws.Range("I31").Value = 25001000
$ws.Range("J31").Value = 7023
$ws.Range("K31").Value = 25001000
$ws.Range("L31").Value = 7023
$ws.Range("M31").Value = -25000752
$ws.Range("N31").Value = -7519

# LTW row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5299.0586
$ws.Range("I40").Value = 3801.4285
$ws.Range("J40").Value = 6347.4
$ws.Range("K40").Value = 3801.4285
$ws.Range("L40").Value = 6347.4
$ws.Range("M40").Value = -3665.4285
$ws.Range("N40").Value = -6619.4

# LTW row 126 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6208.684
$ws.Range("I126").Value = 6166.5
$ws.Range("J126").Value = 6219.933
$ws.Range("K126").Value = 18499.5
$ws.Range("L126").Value = 18659.799
$ws.Range("M126").Value = -16029.5
$ws.Range("N126").Value = -23599.799

# WVR row 81 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 923.4
$ws.Range("I81").Value = 923.4
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1846.8
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -785.8

# WVR row 84 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 923.4
$ws.Range("I84").Value = 923.4
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9234
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -3930

# WVR row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3393.2703
$ws.Range("I132").Value = 2038.0714
$ws.Range("J132").Value = 7609.4443
$ws.Range("K132").Value = 6114.2142
$ws.Range("L132").Value = 22828.3329
$ws.Range("M132").Value = -3584.2142
$ws.Range("N132").Value = -27888.3329

Write-Host "Applied scheduled runner price updates across 8 sheets (36 rows)."
